$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.39), maa://25390 (96.09), maa://36681 (87.34)'
$ws.Range('L2').Value = '*maa://24633 (56.33), *maa://30515 (69.61), *maa://34787 (72.97), ***maa://20792 (11.93), maa://39402 (90.38), ***maa://29083 (27.78)'
$ws.Range('P3').Value = 'maa://21249 (94.3), maa://26254 (96.3)'
$ws.Range('S3').NumberFormat = '@'
$ws.Range('S3').Value = '4'
$ws.Range('T3').Value = 'maa://24617 (89.38), **maa://20790 (44.12), ***maa://37170 (17.19), maa://45854 (100.0)'
$ws.Range('AB3').Value = 'maa://24390 (94.03)'
$ws.Range('D4').Value = 'maa://24632 (93.63), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range('T4').Value = 'maa://32509 (97.3), maa://27295 (84.85), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)'
$ws.Range('X4').Value = '**maa://32495 (48.31), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (85.37)'
$ws.Range('D5').Value = 'maa://21245 (84.0), maa://22744 (84.0)'
$ws.Range('P6').Value = 'maa://31836 (92.31), maa://30381 (92.31)'
$ws.Range('AB6').Value = 'maa://22739 (92.98)'
$ws.Range('A8').Value = '更新日期：2025.01.26 08:43:45'
$ws.Range('X8').Value = 'maa://21411 (96.06)'
$ws.Range('AB8').Value = 'maa://25389 (87.88)'
$ws.Range('AF8').Value = '*maa://24479 (77.91), *maa://21990 (51.85)'
$ws.Range('D9').Value = 'maa://22765 (92.31), *maa://21915 (69.23)'
$ws.Range('P9').Value = 'maa://22736 (82.29)'
$ws.Range('T9').Value = '**maa://22866 (30.19), maa://26222 (97.92)'
$ws.Range('X9').Value = 'maa://26223 (97.67)'
$ws.Range('AF9').Value = 'maa://26206 (90.18), *maa://22865 (51.92)'
$ws.Range('D10').Value = '***maa://25695 (19.02), **maa://32237 (40.91), ***maa://34206 (20.83), ***maa://39951 (15.56), ***maa://39243 (28.57), *maa://45271 (53.33)'
$ws.Range('T10').Value = 'maa://27395 (96.24), maa://22755 (87.72), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('W10').NumberFormat = '@'
$ws.Range('W10').Value = '3'
$ws.Range('X10').Value = 'maa://22301 (97.7), maa://22726 (100.0), maa://45828 (100.0)'
$ws.Range('AF10').Value = '*maa://25021 (54.65), *maa://22733 (60.61), maa://22761 (100.0)'
$ws.Range('AF11').Value = 'maa://31203 (95.83)'
$ws.Range('X12').Value = 'maa://22753 (91.18), *maa://21485 (76.81), maa://37962 (87.88)'
$ws.Range('AB12').Value = 'maa://23669 (95.44), maa://36677 (92.59), maa://39872 (90.91)'
$ws.Range('AF12').Value = '*maa://28932 (78.1), *maa://20106 (63.96), *maa://22769 (64.29)'
$ws.Range('D13').Value = 'maa://24999 (91.82), maa://36673 (92.96), maa://25001 (85.51)'
$ws.Range('P13').Value = 'maa://22676 (92.37), *maa://22583 (74.24), *maa://22500 (57.78)'
$ws.Range('X13').Value = 'maa://34957 (81.43), *maa://22768 (51.61)'
$ws.Range('D14').Value = 'maa://30764 (88.68)'
$ws.Range('P14').Value = 'maa://23250 (98.68), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('D15').Value = '*maa://22743 (77.39), maa://22734 (84.03), *maa://30808 (64.18), **maa://36048 (40.74), maa://45058 (100.0)'
$ws.Range('AF15').Value = 'maa://21364 (81.45), *maa://36666 (78.89), *maa://22766 (69.91)'
$ws.Range('D16').Value = 'maa://21441 (96.35), maa://36679 (93.62), maa://37650 (96.97)'
$ws.Range('X16').Value = 'maa://28501 (97.92), maa://28051 (96.0)'
$ws.Range('AB16').Value = 'maa://26228 (95.7)'
$ws.Range('H17').Value = 'maa://22430 (88.6), maa://39599 (84.44)'
$ws.Range('D18').Value = 'maa://24570 (97.21)'
$ws.Range('AB19').Value = '*maa://30709 (64.54), *maa://36668 (57.5)'
$ws.Range('L20').Value = 'maa://41331 (86.29)'
$ws.Range('T20').Value = 'maa://29113 (85.71)'
$ws.Range('AF21').Value = 'maa://22524 (94.29), *maa://22432 (78.46)'
$ws.Range('L22').Value = 'maa://27127 (82.14), *maa://22751 (71.64)'
$ws.Range('X22').Value = 'maa://21282 (98.54), *maa://37649 (67.86)'
$ws.Range('L23').Value = 'maa://39756 (95.0), maa://39875 (93.94)'
$ws.Range('AB23').Value = 'maa://29652 (97.62)'
$ws.Range('D24').Value = '*maa://24368 (78.11)'
$ws.Range('X24').Value = 'maa://29988 (85.19), maa://23504 (93.09), **maa://22892 (40.14), *maa://25141 (76.74), *maa://36663 (77.78), ***maa://22815 (23.08)'
$ws.Range('AE24').NumberFormat = '@'
$ws.Range('AE24').Value = '5'
$ws.Range('AF24').Value = 'maa://22523 (85.71), maa://36672 (80.36), maa://29910 (92.86), **maa://21440 (34.55), maa://45831 (100.0)'
$ws.Range('L25').Value = 'maa://24378 (87.8)'
$ws.Range('AB25').Value = 'maa://31215 (86.92), *maa://24516 (80.0), maa://26001 (87.5)'
$ws.Range('AB26').Value = 'maa://42235 (94.25)'
$ws.Range('H27').Value = '**maa://21283 (48.0), *maa://39601 (76.47), maa://34494 (97.06), **maa://36665 (44.44)'
$ws.Range('T27').Value = '*maa://30624 (77.19)'
$ws.Range('D28').Value = 'maa://24465 (90.93), maa://25725 (83.72)'
$ws.Range('L28').Value = 'maa://30770 (80.43)'
$ws.Range('X28').Value = 'maa://39929 (90.68), maa://41749 (91.3), ***maa://39723 (14.29)'
$ws.Range('AF28').Value = 'maa://36660 (92.38), *maa://36701 (65.52)'
$ws.Range('H29').Value = '*maa://25175 (66.0)'
$ws.Range('L29').Value = 'maa://28432 (93.17), *maa://28440 (78.85), maa://31400 (100.0), *maa://28650 (71.43)'
$ws.Range('AF29').Value = '*maa://24080 (68.95), maa://42865 (82.22), ***maa://34960 (8.33)'
$ws.Range('AA30').NumberFormat = '@'
$ws.Range('AA30').Value = '3'
$ws.Range('AB30').Value = 'maa://42979 (96.21), maa://45045 (100.0), maa://45822 (100.0)'
$ws.Range('L31').Value = 'maa://35926 (93.48), maa://36258 (84.76), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.45), maa://36667 (98.57), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (96.0), maa://41108 (88.0), maa://41238 (96.88), maa://45523 (100.0)'
$ws.Range('P33').Value = 'maa://21956 (80.28), *maa://22730 (79.31)'
$ws.Range('L35').Value = 'maa://41296 (95.77)'
$ws.Range('P37').Value = 'maa://21280 (89.0), *maa://21239 (66.67)'
$ws.Range('P38').Value = '*maa://24383 (69.0)'
$ws.Range('H39').Value = 'maa://36670 (88.04), maa://25199 (84.82), maa://30434 (90.41), ***maa://25036 (16.0), *maa://45059 (66.67), *maa://44165 (66.67)'
$ws.Range('P39').Value = 'maa://24709 (92.09)'
$ws.Range('S39').NumberFormat = '@'
$ws.Range('S39').Value = '2'
$ws.Range('T39').Value = '*maa://45788 (76.0), *maa://45790 (80.0)'
$ws.Range('H43').Value = 'maa://22525 (92.81), maa://21284 (85.11)'
$ws.Range('H44').Value = 'maa://29768 (97.92), maa://27728 (96.0)'
$ws.Range('H47').Value = 'maa://27410 (96.35), maa://29661 (97.24), maa://28038 (84.62)'
$ws.Range('H53').Value = 'maa://32534 (94.01), **maa://32434 (34.78)'
$ws.Range('H57').Value = 'maa://25176 (98.31)'
